# Apply updates to the "Recorded By" column (G) on the
# "Session Analysis Results" worksheet: reorder the comma-separated
# recorder names for a specific set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows where "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$dnasrRows = @(10,11,12,13,14,15,17,18,19,20,21,22,24,26,36,37,38,39,40,41,43,44,45,46,47,48,50,52,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,136,137,138,142,144,145,146,148,151,153)

foreach ($r in $dnasrRows) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}

# Rows where "backup@backdoor.com, System, system" -> "system, backup@backdoor.com, System"
$backdoorRows = @(2,28,54)

foreach ($r in $backdoorRows) {
    $ws.Cells.Item($r, 7).Value = "system, backup@backdoor.com, System"
}
